$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.389.32"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").Value = "2.252.58"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.14"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.41%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -3.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.69"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.22"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("E13").Value = "  -2.23%  "

$ws.Range("D14").Value = "2.591.83"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.58"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.858"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").Value = "2.253.25"
$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").Value = "42.215.88"
$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("E19").Value = "  +0.90%  "

$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("E22").Value = "  +3.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +36.20%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.51"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.86"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.72"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0826"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.98%  "

$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.16%  "

$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.23"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.45"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0317"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.85"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.92%  "

$ws.Range("E40").Value = "  -3.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "62.88"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.35%  "

$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.99"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.75"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("E49").Value = "  -2.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.33"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.26"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.51%  "
